$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: only B changes
$ws.Range("B21").Value = 91771

# Rows 22 and 23 swap their full content (A, B, D, E, F, G, H, I, Q, R),
# with B also incremented by 4 in each case.
# Save old row 22 values before overwriting (use Value2 - plain read/write,
# unlike Value which is not reliable for reads in this runtime).
$a22 = $ws.Range("A22").Value2
$d22 = $ws.Range("D22").Value2
$e22 = $ws.Range("E22").Value2
$f22 = $ws.Range("F22").Value2
$g22 = $ws.Range("G22").Value2
$h22 = $ws.Range("H22").Value2
$i22 = $ws.Range("I22").Value2
$q22 = $ws.Range("Q22").Value2
$r22 = $ws.Range("R22").Value2

$a23 = $ws.Range("A23").Value2
$d23 = $ws.Range("D23").Value2
$e23 = $ws.Range("E23").Value2
$f23 = $ws.Range("F23").Value2
$g23 = $ws.Range("G23").Value2
$h23 = $ws.Range("H23").Value2
$i23 = $ws.Range("I23").Value2
$q23 = $ws.Range("Q23").Value2
$r23 = $ws.Range("R23").Value2

# New row 22 = old row 23 content, B incremented by 4 (57893 -> 57897)
$ws.Range("A22").Value = $a23
$ws.Range("B22").Value = 57897
$ws.Range("D22").Value = $d23
$ws.Range("E22").Value = $e23
$ws.Range("F22").Value = $f23
$ws.Range("G22").Value = $g23
$ws.Range("H22").Value = $h23
$ws.Range("I22").Value = $i23
$ws.Range("Q22").Value = $q23
$ws.Range("R22").Value = $r23

# New row 23 = old row 22 content, B incremented by 4 (91767 -> 91771)
$ws.Range("A23").Value = $a22
$ws.Range("B23").Value = 91771
$ws.Range("D23").Value = $d22
$ws.Range("E23").Value = $e22
$ws.Range("F23").Value = $f22
$ws.Range("G23").Value = $g22
$ws.Range("H23").Value = $h22
$ws.Range("I23").Value = $i22
$ws.Range("Q23").Value = $q22
$ws.Range("R23").Value = $r22

# Row 24: only B changes
$ws.Range("B24").Value = 92530

# Row 26: only B changes
$ws.Range("B26").Value = 97881

# Row 27: only B changes
$ws.Range("B27").Value = 97878
